# Generate Report for Handoff
# A new handoff (b.*.xlf) was generated for b.md in both the zh-cn and
# de-de locales. Reflect the refreshed status/handoff metadata on the
# Overview sheet as well as on each locale's detail sheet, including the
# new "out of date handback" error surfaced for b.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row for b.md (row 3)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-10-20 00:15:24"

# ---------------------------------------------------------------------
# zh-cn sheet - row for b.md (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-10-20 00:15:12"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/878710e2c0cb00f4c30aba152034a89e736789d2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/21e5602842389779217ce4c21094404878f35617/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet - row for b.md (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-10-20 00:15:24"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/878710e2c0cb00f4c30aba152034a89e736789d2/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/21e5602842389779217ce4c21094404878f35617/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667
